$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content edits (rows 7-9, column A) ---
# Row 7: was the text "data" -> becomes the number 1
$ws.Range("A7").Value = 1
# Row 8: was the number 1 -> becomes the text "xsd:string"
$ws.Range("A8").Value = "xsd:string"
# Row 9: was the text "xsd:string" -> becomes the text "data"
$ws.Range("A9").Value = "data"

# Row 9 used to carry the same ("applyFont") style as row 8; now it should
# match row 7's plain/unstyled look instead.
$ws.Range("A7").Copy()
$ws.Range("A9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Selection / active cell ---
$ws.Range("A8").Select()

# --- Window tab ratio ---
$excel.ActiveWindow.TabRatio = 990
